$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
# A8: "Volume 33   Number  5" -> "Volume 33   Number  6"  (the trailing "5" run becomes "6")
$ws.Range("A8").Characters(21, 1).Text = "6"

# C9: "Report Covering the Week  1/26/2026  Through  2/1/2026"
#   -> "Report Covering the Week  2/2/2026  Through  2/8/2026"
# Replace the later date first so the earlier replacement's differing length
# doesn't shift the character offsets of the one that comes after it.
$ws.Range("C9").Characters(47, 8).Text = "2/8/2026"
$ws.Range("C9").Characters(27, 9).Text = "2/2/2026"

# --- Crime-statistics table updates (rows 14-28) ---

    # Row 14
    $ws.Range("M14").Value = -100
    $ws.Range("M14").NumberFormat = '#,##0.0;"-"#,##0.0'

    # Row 15
    $ws.Range("C15").Value = 1
    $ws.Range("C15").NumberFormat = '#,##0'
    $ws.Range("D15").Value = 1
    $ws.Range("D15").NumberFormat = '#,##0'
    $ws.Range("E15").Value = 0
    $ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("G15").Value = 2
    $ws.Range("H15").Value = 100
    $ws.Range("I15").Value = 5
    $ws.Range("J15").Value = 2
    $ws.Range("K15").Value = 150
    $ws.Range("M15").Value = 400
    $ws.Range("N15").Value = 400

    # Row 16
    $ws.Range("C16").Value = 5
    $ws.Range("D16").Value = 4
    $ws.Range("E16").Value = 25
    $ws.Range("F16").Value = 10
    $ws.Range("G16").Value = 12
    $ws.Range("H16").Value = -16.666666666666
    $ws.Range("I16").Value = 18
    $ws.Range("J16").Value = 15
    $ws.Range("K16").Value = 20
    $ws.Range("L16").Value = -30.769230769230
    $ws.Range("M16").Value = -41.935483870967
    $ws.Range("N16").Value = -88.461538461538

    # Row 17
    $ws.Range("C17").Value = 5
    $ws.Range("D17").Value = 6
    $ws.Range("E17").Value = -16.666666666666
    $ws.Range("F17").Value = 19
    $ws.Range("G17").Value = 18
    $ws.Range("H17").Value = 5.555555555555
    $ws.Range("I17").Value = 25
    $ws.Range("J17").Value = 24
    $ws.Range("K17").Value = 4.166666666666
    $ws.Range("L17").Value = 19.047619047619
    $ws.Range("M17").Value = 92.307692307692
    $ws.Range("N17").Value = -16.666666666666

    # Row 18
    $ws.Range("D18").Value = 8
    $ws.Range("E18").Value = -75
    $ws.Range("F18").Value = 9
    $ws.Range("G18").Value = 22
    $ws.Range("H18").Value = -59.090909090909
    $ws.Range("I18").Value = 18
    $ws.Range("J18").Value = 30
    $ws.Range("K18").Value = -40
    $ws.Range("L18").Value = -25
    $ws.Range("M18").Value = -28
    $ws.Range("N18").Value = -88.535031847133

    # Row 19
    $ws.Range("C19").Value = 20
    $ws.Range("D19").Value = 15
    $ws.Range("E19").Value = 33.333333333333
    $ws.Range("F19").Value = 59
    $ws.Range("G19").Value = 45
    $ws.Range("H19").Value = 31.111111111111
    $ws.Range("I19").Value = 79
    $ws.Range("J19").Value = 54
    $ws.Range("K19").Value = 46.296296296296
    $ws.Range("L19").Value = -7.058823529411
    $ws.Range("M19").Value = 83.720930232558
    $ws.Range("N19").Value = -7.058823529411

    # Row 20
    $ws.Range("C20").Value = 1
    $ws.Range("D20").Value = 2
    $ws.Range("E20").Value = -50
    $ws.Range("F20").Value = 24
    $ws.Range("G20").Value = 8
    $ws.Range("H20").Value = 200
    $ws.Range("I20").Value = 35
    $ws.Range("J20").Value = 13
    $ws.Range("K20").Value = 169.230769230769
    $ws.Range("L20").Value = 169.230769230769
    $ws.Range("M20").Value = 40
    $ws.Range("N20").Value = -84.848484848484

    # Row 21
    $ws.Range("C21").Value = 34
    $ws.Range("D21").Value = 36
    $ws.Range("E21").Value = -5.555555555555
    $ws.Range("F21").Value = 125
    $ws.Range("G21").Value = 108
    $ws.Range("H21").Value = 15.740740740740
    $ws.Range("I21").Value = 180
    $ws.Range("J21").Value = 139
    $ws.Range("K21").Value = 29.496402877697
    $ws.Range("L21").Value = 6.508875739644
    $ws.Range("M21").Value = 29.496402877697
    $ws.Range("N21").Value = -72.768532526475

    # Row 22
    $ws.Range("C22").Value = 4
    $ws.Range("F22").Value = 8
    $ws.Range("H22").Value = 700
    $ws.Range("I22").Value = 10
    $ws.Range("K22").Value = 900
    $ws.Range("L22").Value = 233.333333333333
    $ws.Range("M22").Value = 100

    # Row 24
    $ws.Range("C24").Value = 29
    $ws.Range("D24").Value = 30
    $ws.Range("E24").Value = -3.333333333333
    $ws.Range("F24").Value = 116
    $ws.Range("G24").Value = 118
    $ws.Range("H24").Value = -1.694915254237
    $ws.Range("I24").Value = 153
    $ws.Range("J24").Value = 179
    $ws.Range("K24").Value = -14.525139664804
    $ws.Range("L24").Value = -30.769230769230
    $ws.Range("M24").Value = 106.756756756757

    # Row 25
    $ws.Range("C25").Value = 17
    $ws.Range("D25").Value = 12
    $ws.Range("E25").Value = 41.666666666666
    $ws.Range("G25").Value = 65
    $ws.Range("H25").Value = -1.538461538461
    $ws.Range("I25").Value = 84
    $ws.Range("J25").Value = 103
    $ws.Range("K25").Value = -18.446601941747
    $ws.Range("L25").Value = -36.842105263157

    # Row 26
    $ws.Range("C26").Value = 4
    $ws.Range("D26").Value = 4
    $ws.Range("E26").Value = 0
    $ws.Range("F26").Value = 35
    $ws.Range("G26").Value = 26
    $ws.Range("H26").Value = 34.615384615384
    $ws.Range("I26").Value = 47
    $ws.Range("J26").Value = 42
    $ws.Range("K26").Value = 11.904761904761
    $ws.Range("L26").Value = -6
    $ws.Range("M26").Value = -14.545454545454

    # Row 27
    $ws.Range("C27").Value = 2
    $ws.Range("C27").NumberFormat = '#,##0'
    $ws.Range("D27").Value = 2
    $ws.Range("D27").NumberFormat = '#,##0'
    $ws.Range("E27").Value = 0
    $ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("F27").Value = 5
    $ws.Range("G27").Value = 3
    $ws.Range("H27").Value = 66.666666666666
    $ws.Range("I27").Value = 6
    $ws.Range("J27").Value = 4
    $ws.Range("K27").Value = 50
    $ws.Range("L27").Value = 500

    # Row 28
    $ws.Range("D28").Value = 2
    $ws.Range("D28").NumberFormat = '#,##0'
    $ws.Range("E28").Value = 0
    $ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
    $ws.Range("F28").Value = 4
    $ws.Range("G28").Value = 2
    $ws.Range("I28").Value = 5
    $ws.Range("J28").Value = 3
    $ws.Range("K28").Value = 66.666666666666
    $ws.Range("L28").Value = -28.571428571428

